$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "寻址方式（可选）" (Addressing method - optional) column in column O,
# with "静态" (Static) / "动态" (Dynamic) sample values in rows 2 and 3.
$ws.Range("O1").Value = "寻址方式（可选）"
$ws.Range("O2").Value = "静态"
$ws.Range("O3").Value = "动态"

# Match the new column width recorded in the target workbook (~18.875 characters).
$ws.Range("O1").EntireColumn.ColumnWidth = 18.14

# Update the active selection to O3, matching the saved workbook view state.
$null = $ws.Range("O3").Select()
